# "Generate Report for Handoff"
# Update the localization-status report: move Status from
# "Handed back: in sync with en-US" to "Ready for handoff" and refresh the
# generation timestamps, then resize the now-shorter Status columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value = "Ready for handoff"       # zh-cn sheet Status column
$wsDeDe.Range("C2").Value = "Ready for handoff"       # de-de sheet Status column

# --- Refresh generation timestamps ---
# Overview "Latest HO Xliff Generate Date" and de-de "Latest Handoff Datetime"
# shared the same timestamp; bump both to the new handoff time.
$wsOverview.Range("G2").Value = "2016-08-20 19:07:28"
$wsDeDe.Range("H2").Value = "2016-08-20 19:07:28"
# zh-cn "Latest Handoff Datetime" gets its own refreshed timestamp.
$wsZhCn.Range("H2").Value = "2016-08-20 19:07:24"

# --- Resize the Status columns now that the text is shorter ---
$wsOverview.Columns(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns(3).ColumnWidth = 16.333333333333332
